$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D values are text (prices formatted as strings in the source
# data, some containing multiple "." separators or subscript digits). A
# leading apostrophe is used for numeric-looking values so Excel stores them
# as text (matching the original inlineStr/text cell type) instead of
# auto-converting them to numbers and silently dropping trailing zeros.

$ws.Range("D2").Value = "54.148.84"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "2.270.35"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").Value = "'497.76"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").Value = "'128.12"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("E11").Value = "  +3.40%  "

$ws.Range("D12").Value = "'4.71"
$ws.Range("E12").Value = "  +1.90%  "

$ws.Range("D13").Value = "2.670.88"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").Value = "'22.62"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").Value = "54.124.88"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "2.266.66"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").Value = "'10.19"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").Value = "'302.78"
$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("D23").Value = "'61.04"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("E26").Value = "  +2.32%  "

$ws.Range("D27").Value = "'170.71"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").Value = "0.0₃0689"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'5.92"
$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "'0.943"
$ws.Range("E35").Value = "  +9.39%  "

$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").Value = "'3.69"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("E38").Value = "  -1.21%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D44").Value = "'0.0890"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").Value = "'238.78"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").Value = "'0.0205"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'16.15"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("E51").Value = "  -0.46%  "

# Row 41/42 swap: Aave moves to row 41, RenderToken moves to row 42
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'124.86"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'4.80"
$ws.Range("E42").Value = "  -2.39%  "
